# Checkout: deduct ingredients used for the order from inventory, and
# reflect the resulting totals in row 2 of the inventory sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 988
$ws.Range("B2").Value = 788
$ws.Range("C2").Value = 788
$ws.Range("D2").Value = 788
$ws.Range("G2").Value = 988
